# Updates cryptos list values (Price + Volume(1h)) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.242.49"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "1.566.05"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'210.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'0.0872"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "1.788.30"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "1.520.04"
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "27.195.02"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "'61.96"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'217.96"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "'151.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "'6.63"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").Value = "1.458.57"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").Value = "'3.17"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  +3.97%  "
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'2.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").Value = "'64.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "1.699.43"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'85.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'0.0526"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").Value = "'0.0948"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.18%  "
